$wb = $excel.ActiveWorkbook

# Duplicate the "Croatia" sheet (closest template for the new market sheet),
# placing the copy right after it, then rename it to "Greece".
$croatia = $wb.Worksheets.Item("Croatia")
$croatia.Copy($null, $croatia)
$greece = $wb.Worksheets.Item($wb.Worksheets.Count)
$greece.Name = "Greece"

# Update the market-specific content on the new sheet.
$greece.Range("B2").Value = "Greece Market"
$greece.Range("B4").Value = "NGC-4119/T3164"

# Column widths differ slightly from the Croatia template.
$greece.Columns.Item(1).ColumnWidth = 23.608072916666668
$greece.Columns.Item(2).ColumnWidth = 16.944010416666668

# Rows 3-5 drop their custom (taller) row height back to the sheet default.
$greece.Rows.Item(3).EntireRow.AutoFit()
$greece.Rows.Item(4).EntireRow.AutoFit()
$greece.Rows.Item(5).EntireRow.AutoFit()

# Croatia's remembered selection moves on (it is no longer the active tab).
$croatia.Activate()
$croatia.Range("I23").Select()

# Greece becomes the active/selected tab with its own remembered selection.
$greece.Activate()
$greece.Range("C15").Select()
